$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$origStyle_D2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.059.57"
$ws.Range("D2").Style = $origStyle_D2
$ws.Range("E2").Value = "  +0.37%  "

# Row 3
$origStyle_D3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.363.21"
$ws.Range("D3").Style = $origStyle_D3
$ws.Range("E3").Value = "  +0.16%  "

# Row 4
$ws.Range("E4").Value = "  +0.25%  "

# Row 5
$origStyle_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.694"
$ws.Range("D5").Style = $origStyle_D5
$ws.Range("E5").Value = "  +5.16%  "

# Row 6
$origStyle_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.90"
$ws.Range("D6").Style = $origStyle_D6
$ws.Range("E6").Value = "  +2.43%  "

# Row 7
$origStyle_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.09"
$ws.Range("D7").Style = $origStyle_D7
$ws.Range("E7").Value = "  +3.67%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$origStyle_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.613"
$ws.Range("D9").Style = $origStyle_D9
$ws.Range("E9").Value = "  +14.97%  "

# Row 10
$origStyle_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("D10").Style = $origStyle_D10
$ws.Range("E10").Value = "  +2.45%  "

# Row 11
$origStyle_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.57"
$ws.Range("D11").Style = $origStyle_D11
$ws.Range("E11").Value = "  +1.14%  "

# Row 12
$origStyle_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.31"
$ws.Range("D12").Style = $origStyle_D12
$ws.Range("E12").Value = "  +16.90%  "

# Row 13
$origStyle_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.44"
$ws.Range("D13").Style = $origStyle_D13
$ws.Range("E13").Value = "  +12.07%  "

# Row 14
$ws.Range("E14").Value = "  +2.00%  "

# Row 15
$origStyle_D15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.718.02"
$ws.Range("D15").Style = $origStyle_D15
$ws.Range("E15").Value = "  +0.31%  "

# Row 16
$origStyle_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.65"
$ws.Range("D16").Style = $origStyle_D16
$ws.Range("E16").Value = "  -1.14%  "

# Row 17
$origStyle_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.923"
$ws.Range("D17").Style = $origStyle_D17
$ws.Range("E17").Value = "  +4.30%  "

# Row 18
$origStyle_D18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.361.49"
$ws.Range("D18").Style = $origStyle_D18
$ws.Range("E18").Value = "  +0.25%  "

# Row 19
$origStyle_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.954.00"
$ws.Range("D19").Style = $origStyle_D19
$ws.Range("E19").Value = "  +0.51%  "

# Row 20
$origStyle_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000102"
$ws.Range("D20").Style = $origStyle_D20
$ws.Range("E20").Value = "  +1.45%  "

# Row 21
$origStyle_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.71"
$ws.Range("D21").Style = $origStyle_D21
$ws.Range("E21").Value = "  +6.23%  "

# Row 22
$origStyle_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.52"
$ws.Range("D22").Style = $origStyle_D22
$ws.Range("E22").Value = "  +1.89%  "

# Row 23
$origStyle_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "259.01"
$ws.Range("D23").Style = $origStyle_D23
$ws.Range("E23").Value = "  +3.21%  "

# Row 24
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
$ws.Range("E25").Value = "  -1.55%  "

# Row 26
$origStyle_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.52"
$ws.Range("D26").Style = $origStyle_D26
$ws.Range("E26").Value = "  +1.56%  "

# Row 27
$origStyle_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.79"
$ws.Range("D27").Style = $origStyle_D27
$ws.Range("E27").Value = "  +16.23%  "

# Row 28
$origStyle_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.85"
$ws.Range("D28").Style = $origStyle_D28
$ws.Range("E28").Value = "  +5.55%  "

# Row 29
$ws.Range("E29").Value = "  +1.51%  "

# Row 30
$origStyle_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.15"
$ws.Range("D30").Style = $origStyle_D30
$ws.Range("E30").Value = "  +2.72%  "

# Row 31
$origStyle_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.03"
$ws.Range("D31").Style = $origStyle_D31
$ws.Range("E31").Value = "  +1.24%  "

# Row 32
$origStyle_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.127"
$ws.Range("D32").Style = $origStyle_D32
$ws.Range("E32").Value = "  -3.72%  "

# Row 33
$origStyle_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.136"
$ws.Range("D33").Style = $origStyle_D33
$ws.Range("E33").Value = "  +4.36%  "

# Row 34
$origStyle_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.39"
$ws.Range("D34").Style = $origStyle_D34
$ws.Range("E34").Value = "  +4.44%  "

# Row 35
$origStyle_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0762"
$ws.Range("D35").Style = $origStyle_D35
$ws.Range("E35").Value = "  +7.93%  "

# Row 36
$origStyle_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.40"
$ws.Range("D36").Style = $origStyle_D36
$ws.Range("E36").Value = "  +5.03%  "

# Row 37
$origStyle_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.77"
$ws.Range("D37").Style = $origStyle_D37
$ws.Range("E37").Value = "  +0.33%  "

# Row 38
$origStyle_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.42"
$ws.Range("D38").Style = $origStyle_D38
$ws.Range("E38").Value = "  -1.09%  "

# Row 39
$origStyle_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.42"
$ws.Range("D39").Style = $origStyle_D39
$ws.Range("E39").Value = "  -0.05%  "

# Row 40
$origStyle_D40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0282"
$ws.Range("D40").Style = $origStyle_D40
$ws.Range("E40").Value = "  +6.63%  "

# Row 41
$origStyle_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.218"
$ws.Range("D41").Style = $origStyle_D41
$ws.Range("E41").Value = "  +21.29%  "

# Row 42
$origStyle_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.19"
$ws.Range("D42").Style = $origStyle_D42
$ws.Range("E42").Value = "  +3.28%  "

# Row 43
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$origStyle_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.108"
$ws.Range("D43").Style = $origStyle_D43
$ws.Range("E43").Value = "  +12.11%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$origStyle_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.20"
$ws.Range("D44").Style = $origStyle_D44
$ws.Range("E44").Value = "  -1.16%  "

# Row 45
$origStyle_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.88"
$ws.Range("D45").Style = $origStyle_D45
$ws.Range("E45").Value = "  +9.77%  "

# Row 46
$ws.Range("E46").Value = "  +0.04%  "

# Row 47
$origStyle_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.54"
$ws.Range("D47").Style = $origStyle_D47
$ws.Range("E47").Value = "  +10.66%  "

# Row 48
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$origStyle_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.26"
$ws.Range("D48").Style = $origStyle_D48
$ws.Range("E48").Value = "  +3.33%  "

# Row 49
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$origStyle_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.16"
$ws.Range("D49").Style = $origStyle_D49
$ws.Range("E49").Value = "  +3.45%  "

# Row 50
$origStyle_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.19"
$ws.Range("D50").Style = $origStyle_D50
$ws.Range("E50").Value = "  +1.38%  "

# Row 51
$origStyle_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.90"
$ws.Range("D51").Style = $origStyle_D51
$ws.Range("E51").Value = "  +8.01%  "
